# Add 2022-Q1 data:
#  - insert a new worksheet "2022-Q1" (holding the per-fund breakdown) right
#    before the existing "总计" (totals) sheet
#  - prepend a new "2022-Q1" summary row to the "总计" sheet

$wb = $excel.ActiveWorkbook

# Helper: write a value as literal text (keeps leading zeros / avoids
# Excel's automatic number coercion) using the classic leading-quote text
# marker. Only needed for values that "look like" numbers.
function Set-TextCell($cell, $val) {
    $cell.Value = "'" + $val
}

# ---------------------------------------------------------------------
# 1. Locate the current "总计" sheet (3rd tab) and insert a brand new
#    sheet immediately before it, then rename the new sheet "2022-Q1".
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item(3)
$q1Sheet = $wb.Worksheets.Add($totalSheet)
$q1Sheet.Name = "2022-Q1"

# After insertion the old totals sheet has shifted one slot to the right.
$totalSheet = $wb.Worksheets.Item(4)

# ---------------------------------------------------------------------
# 2. Populate the "2022-Q1" sheet.
#    Re-use the existing formatting from the "总计" sheet: the bold
#    bordered header style lives on row 1 / column A of that sheet.
#    Apply the style (via Copy) BEFORE writing any quote-prefixed text so
#    the destination cell never picks up a stray quote-prefix flag.
# ---------------------------------------------------------------------
$headerStyleCell = $totalSheet.Cells.Item(1, 2)
$colAStyleCell = $totalSheet.Cells.Item(2, 1)

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $cell = $q1Sheet.Cells.Item(1, $c + 2)
    $headerStyleCell.Copy($cell)
    # Header text is non-numeric, so a plain assignment already stays text.
    $cell.Value = $headers[$c]
}

$rows = @(
    @("110005", "易方达积极成长混合",               "36.95", "92.57", "3.66", "1.3524", 9),
    @("001725", "汇添富中国高端制造股票",             "18.26", "90.98", "6.33", "1.1559", 1),
    @("015115", "汇添富中国高端制造股票D",            "18.26", "90.98", "6.33", "1.1559", 1),
    @("007639", "汇添富3年封闭运作竞争优势灵活配置混合", "13.07", "60.50", "4.63", "0.6051", 3),
    @("010518", "华夏先锋科技一年定期开放混合A",        "8.72",  "88.08", "3.90", "0.3401", 10),
    @("000690", "前海开源大海洋战略经济灵活配置混合",    "6.91",  "93.87", "4.02", "0.2778", 9),
    @("010412", "汇安均衡优选混合",                   "6.17",  "94.73", "4.08", "0.2517", 10),
    @("002746", "汇添富多策略定期开放灵活配置混合",      "4.55",  "64.75", "4.54", "0.2066", 3),
    @("002181", "华安大安全主题灵活配置混合",           "5.02",  "87.28", "3.06", "0.1536", 8),
    @("011506", "建信高端装备股票型证券投资基金A",       "2.59",  "85.91", "4.44", "0.1150", 3),
    @("501063", "汇添富悦享定期开放混合",               "2.18",  "60.21", "4.41", "0.0961", 4),
    @("010519", "华夏先锋科技一年定期开放混合C",        "2.12",  "88.08", "3.90", "0.0827", 10),
    @("011507", "建信高端装备股票型证券投资基金C",       "0.90",  "85.91", "4.44", "0.0400", 3),
    @("009317", "金信核心竞争力灵活配置混合",            "0.19",  "89.48", "2.94", "0.0056", 8)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = $r + 2

    $aCell = $q1Sheet.Cells.Item($excelRow, 1)
    $colAStyleCell.Copy($aCell)
    $aCell.Value = $r

    # Fund name (column C) is plain non-numeric text - no quoting needed.
    $q1Sheet.Cells.Item($excelRow, 3).Value = $row[1]
    # Fund code + numeric-looking metrics must stay text (leading zeros /
    # trailing zeros would otherwise be stripped by Excel's auto-coercion).
    Set-TextCell $q1Sheet.Cells.Item($excelRow, 2) $row[0]
    Set-TextCell $q1Sheet.Cells.Item($excelRow, 4) $row[2]
    Set-TextCell $q1Sheet.Cells.Item($excelRow, 5) $row[3]
    Set-TextCell $q1Sheet.Cells.Item($excelRow, 6) $row[4]
    Set-TextCell $q1Sheet.Cells.Item($excelRow, 7) $row[5]
    $q1Sheet.Cells.Item($excelRow, 8).Value = $row[6]
}

[void]$q1Sheet.Range("A1").Select()

# ---------------------------------------------------------------------
# 3. Prepend the 2022-Q1 summary row to the "总计" sheet, pushing the
#    existing rows down by one.
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()
# Inserting a row clones formatting from the row above (the bold header);
# the new data row itself carries no special formatting in the target, so
# strip it back to the default first.
$totalSheet.Rows.Item(2).ClearFormats()

$newA = $totalSheet.Cells.Item(2, 1)
$totalSheet.Cells.Item(3, 1).Copy($newA)
$newA.Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 14
$totalSheet.Cells.Item(2, 4).Value = 5.84

$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2

[void]$totalSheet.Range("A1").Select()
